# Fix the "Yapacany" -> "Yapacani" typo in the "S.J. Yapacany" / "San Juan de
# Yapacany" location entries (columns F/G "Siglas"/"Significado" lookup table
# and the D21 data row that uses the same abbreviation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "S.J. Yapacani"
$ws.Range("G4").Value = "San Juan de Yapacani"
$ws.Range("D21").Value = "S.J. Yapacani"

# Update the window selection to match the new active cell/view reflected in
# the saved workbook.
$ws.Range("D24").Select()
